$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Insert a new row above row 10; this pushes the old rows 10-14 down to 11-15
$ws.Rows.Item(10).Insert()

# Fill the new row 10 with its data
$ws.Cells.Item(10, 1).Value = "Problem:Is problem more frequent or once in awhile when the engine is warmed up?"
$ws.Cells.Item(10, 2).Value = "Possible_Problem"
$ws.Cells.Item(10, 3).Value = "Possible_Problem:50% Oil Pressure Sensor / Switch`n25% Oil Pressure Sensor Wiring / Connector`n10% Oil Pressure Control Valve`n5% Internal Engine Failure`n5% Engine Oil Pump`n5% Engine Control Module"

# Row height for the new row matches the target ht="388.8"
$ws.Rows.Item(10).RowHeight = 388.8

# Update view state to match target (scroll to show row 10, select C10)
$ws.Range("C10").Select()
